# 3-c-1.xlsx update: indicator description shortened, organisation website
# domain corrected, and final selection left on the edited "website" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Пример")

# B4 ("Индикатор"): drop the two enumerated sub-bullets, keep only the
# heading line (with its trailing newline).
$ws.Range("B4").Value = "3.c.1. Число медицинских работников на душу населения и их распределение`n"

# B10 ("Сайт организации"): the statistics committee's domain changed.
$ws.Range("B10").Value = "www.stat.gov.kg"

# Leave the cursor on the cell that was last edited, matching the saved
# selection state.
$ws.Range("B10").Select()

Write-Output "Updated B4 and B10 on sheet '$($ws.Name)'"
